# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (spreadsheet column G, header label "K" in G1) was
# recalculated with the corrected (K, not Strike#) values. Rewrite the
# recalculated column G values for rows 2-66 in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 2
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 0
    20 = 0
    21 = 1
    22 = 3
    23 = 0
    24 = 0
    25 = 1
    26 = 2
    27 = 1
    28 = 0
    29 = 0
    30 = 0
    31 = 1
    32 = 1
    33 = 1
    34 = 0
    35 = 3
    36 = 1
    37 = 2
    38 = 2
    39 = 1
    40 = 1
    41 = 2
    42 = 1
    43 = 1
    44 = 0
    45 = 0
    46 = 1
    47 = 1
    48 = 2
    49 = 1
    50 = 1
    51 = 0
    52 = 1
    53 = 0
    54 = 2
    55 = 1
    56 = 1
    57 = 0
    58 = 1
    59 = 1
    60 = 0
    61 = 0
    62 = 2
    63 = 0
    64 = 2
    65 = 1
    66 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
